$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TestCases")
$ws2 = $wb.Worksheets.Item("TestData")

# ---------------------------------------------------------------------------
# TestData sheet: add a new "browser" column to both data-driven test blocks,
# add a new row of test data to the AddCustomerTest block, rewrite the
# OpenAccountTest block with new customer names / currencies, and flip the
# remaining Runmode flags from N to Y.
# ---------------------------------------------------------------------------

# AddCustomerTest block (rows 2-5): new "browser" column E
$ws2.Cells.Item(2, 5).Value = "browser"
$ws2.Cells.Item(3, 5).Value = "chrome"

$ws2.Cells.Item(4, 1).Value = "Y"
$ws2.Cells.Item(4, 5).Value = "firefox"

$ws2.Cells.Item(5, 1).Value = "Y"
$ws2.Cells.Item(5, 2).Value = "Tom"
$ws2.Cells.Item(5, 3).Value = "Smith"
$ws2.Cells.Item(5, 4).Value = "D53434"
$ws2.Cells.Item(5, 5).Value = "firefox"

# Old OpenAccountTest header (row 6) moves down - clear the old spot
$ws2.Cells.Item(6, 1).ClearContents()

# OpenAccountTest block now starts at row 7, header at row 8
$ws2.Cells.Item(7, 1).Value = "OpenAccountTest"
$ws2.Cells.Item(7, 2).ClearContents()
$ws2.Cells.Item(7, 3).ClearContents()

$ws2.Cells.Item(8, 1).Value = "Runmode"
$ws2.Cells.Item(8, 2).Value = "customer"
$ws2.Cells.Item(8, 3).Value = "currency"
$ws2.Cells.Item(8, 4).Value = "browser"

$ws2.Cells.Item(9, 2).Value = "Harry Potter"
$ws2.Cells.Item(9, 3).Value = "Rupee"
$ws2.Cells.Item(9, 4).Value = "chrome"

$ws2.Cells.Item(10, 2).Value = "Ron Weasly"
$ws2.Cells.Item(10, 4).Value = "chrome"

$ws2.Cells.Item(11, 1).Value = "Y"
$ws2.Cells.Item(11, 2).Value = "Hermoine Granger"
$ws2.Cells.Item(11, 3).Value = "Pound"
$ws2.Cells.Item(11, 4).Value = "firefox"

# ---------------------------------------------------------------------------
# View state: TestCases is no longer the active tab, its selection moves to
# B3; TestData becomes the active tab with selection on G5.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("B3").Select()

$ws2.Activate()
$ws2.Range("G5").Select()
